$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 27.93616396685348
$ws.Range("C2").Value = 9.124780832205763
$ws.Range("D2").Value = 4.316218585986169
$ws.Range("F2").Value = 65.45344463466692
$ws.Range("G2").Value = 3.803776099455128
$ws.Range("J2").Value = 10.59259957355561
$ws.Range("L2").Value = 12.14206506012712
$ws.Range("B3").Value = 27.81012124152216
$ws.Range("C3").Value = 8.873488018575655
$ws.Range("D3").Value = 4.185567368798859
$ws.Range("F3").Value = 64.70702755567935
$ws.Range("G3").Value = 3.808897594652299
$ws.Range("J3").Value = 10.59529149339129
$ws.Range("L3").Value = 12.18397200724808
$ws.Range("B4").Value = 27.7429978663978
$ws.Range("C4").Value = 8.720932190553775
$ws.Range("D4").Value = 4.103242631457853
$ws.Range("F4").Value = 64.25351847646543
$ws.Range("G4").Value = 3.81219990739955
$ws.Range("J4").Value = 10.59742139426621
$ws.Range("L4").Value = 12.21216942787957
$ws.Range("B5").Value = 27.71824816476131
$ws.Range("C5").Value = 8.659330999635383
$ws.Range("D5").Value = 4.069198488755617
$ws.Range("F5").Value = 64.07004628130055
$ws.Range("G5").Value = 3.81358545961158
$ws.Range("J5").Value = 10.5984090395837
$ws.Range("L5").Value = 12.22428010400713
$ws.Range("B6").Value = 27.71429630038424
$ws.Range("C6").Value = 8.649140223708169
$ws.Range("D6").Value = 4.063516533951025
$ws.Range("F6").Value = 64.0396653320991
$ws.Range("G6").Value = 3.813817940571028
$ws.Range("J6").Value = 10.59858025873673
$ws.Range("L6").Value = 12.22632851308558
$ws.Range("B7").Value = 27.7426535156204
$ws.Range("C7").Value = 8.720098948289007
$ws.Range("D7").Value = 4.102785463616801
$ws.Range("F7").Value = 64.25103852488024
$ws.Range("G7").Value = 3.81221843191948
$ws.Range("J7").Value = 10.59743422971065
$ws.Range("L7").Value = 12.21233024658281
$ws.Range("B8").Value = 27.89058782763919
$ws.Range("C8").Value = 9.037850613507763
$ws.Range("D8").Value = 4.271625813025349
$ws.Range("F8").Value = 65.19516462393614
$ws.Range("G8").Value = 3.805509370012468
$ws.Range("J8").Value = 10.59342857565785
$ws.Range("L8").Value = 12.1560025624465
$ws.Range("B9").Value = 28.2609927406265
$ws.Range("C9").Value = 9.66928789877379
$ws.Range("D9").Value = 4.584707133451374
$ws.Range("F9").Value = 67.07874340178684
$ws.Range("G9").Value = 3.793595950673619
$ws.Range("J9").Value = 10.58937114147143
$ws.Range("L9").Value = 12.06512288132661
$ws.Range("B10").Value = 28.5802526129279
$ws.Range("C10").Value = 10.13136126540055
$ws.Range("D10").Value = 4.802168934642752
$ws.Range("F10").Value = 68.47419865014346
$ws.Range("G10").Value = 3.785589449438366
$ws.Range("J10").Value = 10.5887232676665
$ws.Range("L10").Value = 12.01030223408555
$ws.Range("B11").Value = 28.7352552307332
$ws.Range("C11").Value = 10.33984797199868
$ws.Range("D11").Value = 4.898088323308341
$ws.Range("F11").Value = 69.10982222383942
$ws.Range("G11").Value = 3.782106662376856
$ws.Range("J11").Value = 10.58893903614533
$ws.Range("L11").Value = 11.98796015395249
$ws.Range("B12").Value = 28.79531103269106
$ws.Range("C12").Value = 10.41845111477407
$ws.Range("D12").Value = 4.933957944403954
$ws.Range("F12").Value = 69.35049006566916
$ws.Range("G12").Value = 3.780810552035176
$ws.Range("J12").Value = 10.5890944473739
$ws.Range("L12").Value = 11.97987332852969
$ws.Range("B13").Value = 28.78231721992222
$ws.Range("C13").Value = 10.40153942263282
$ws.Range("D13").Value = 4.926253233443818
$ws.Range("F13").Value = 69.29866156149377
$ws.Range("G13").Value = 3.781088683801099
$ws.Range("J13").Value = 10.58905769406243
$ws.Range("L13").Value = 11.98159834859755
$ws.Range("B14").Value = 28.74016907710896
$ws.Range("C14").Value = 10.34632218164196
$ws.Range("D14").Value = 4.901048528109516
$ws.Range("F14").Value = 69.12962353270878
$ws.Range("G14").Value = 3.781999575762421
$ws.Range("J14").Value = 10.58895034307999
$ws.Range("L14").Value = 11.98728735597372
$ws.Range("B15").Value = 28.71452780071476
$ws.Range("C15").Value = 10.31245211328141
$ws.Range("D15").Value = 4.885550348610819
$ws.Range("F15").Value = 69.02607472989907
$ws.Range("G15").Value = 3.782560480257702
$ws.Range("J15").Value = 10.58889419477089
$ws.Range("L15").Value = 11.99082070483835
$ws.Range("B16").Value = 28.57031562356884
$ws.Range("C16").Value = 10.11769296641164
$ws.Range("D16").Value = 4.795838073073522
$ws.Range("F16").Value = 68.43266309420352
$ws.Range("G16").Value = 3.785820250340016
$ws.Range("J16").Value = 10.58871946443709
$ws.Range("L16").Value = 12.01181460939016
$ws.Range("B17").Value = 28.48431718560718
$ws.Range("C17").Value = 9.997702063925027
$ws.Range("D17").Value = 4.740017263276957
$ws.Range("F17").Value = 68.06872994984374
$ws.Range("G17").Value = 3.787860716306696
$ws.Range("J17").Value = 10.58874323322628
$ws.Range("L17").Value = 12.02535884572398
$ws.Range("B18").Value = 28.43577563715466
$ws.Range("C18").Value = 9.92853070202553
$ws.Range("D18").Value = 4.707628853769751
$ws.Range("F18").Value = 67.85949155117206
$ws.Range("G18").Value = 3.789049353460193
$ws.Range("J18").Value = 10.58880493648406
$ws.Range("L18").Value = 12.03339347197049
$ws.Range("B19").Value = 28.41950006695965
$ws.Range("C19").Value = 9.905086897717393
$ws.Range("D19").Value = 4.696614975792401
$ws.Range("F19").Value = 67.78866636685576
$ws.Range("G19").Value = 3.789454389836797
$ws.Range("J19").Value = 10.58883406864225
$ws.Range("L19").Value = 12.03615581399492
$ws.Range("B20").Value = 28.4933767110677
$ws.Range("C20").Value = 10.01049217265084
$ws.Range("D20").Value = 4.745988802705607
$ws.Range("F20").Value = 68.10746339362638
$ws.Range("G20").Value = 3.787641952501858
$ws.Range("J20").Value = 10.58873572970701
$ws.Range("L20").Value = 12.02389174807073
$ws.Range("B21").Value = 28.7525124818663
$ws.Range("C21").Value = 10.36255097553127
$ws.Range("D21").Value = 4.908464209349709
$ws.Range("F21").Value = 69.17927605429624
$ws.Range("G21").Value = 3.781731408821096
$ws.Range("J21").Value = 10.58897987199505
$ws.Range("L21").Value = 11.98560621377018
$ws.Range("B22").Value = 28.92977495675559
$ws.Range("C22").Value = 10.59058316744923
$ws.Range("D22").Value = 5.012002935173273
$ws.Range("F22").Value = 69.87953135140289
$ws.Range("G22").Value = 3.778001022684562
$ws.Range("J22").Value = 10.58956917531445
$ws.Range("L22").Value = 11.96276237535708
$ws.Range("B23").Value = 28.83445925059422
$ws.Range("C23").Value = 10.46909729062549
$ws.Range("D23").Value = 4.956990996792294
$ws.Range("F23").Value = 69.50586182007207
$ws.Range("G23").Value = 3.779979935291779
$ws.Range("J23").Value = 10.58921523398175
$ws.Range("L23").Value = 11.97475517549612
$ws.Range("B24").Value = 28.48927809260732
$ws.Range("C24").Value = 10.00471033926207
$ws.Range("D24").Value = 4.743289990977026
$ws.Range("F24").Value = 68.08995201861028
$ws.Range("G24").Value = 3.787740807181795
$ws.Range("J24").Value = 10.58873897244208
$ws.Range("L24").Value = 12.02455425090205
$ws.Range("B25").Value = 28.15236496694595
$ws.Range("C25").Value = 9.498316466806047
$ws.Range("D25").Value = 4.50211165504471
$ws.Range("F25").Value = 66.56656485664114
$ws.Range("G25").Value = 3.796686959555316
$ws.Range("J25").Value = 10.59006035461987
$ws.Range("L25").Value = 12.08761083185621